$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "277.88"
Set-TextValue "E2" "0.90%"

Set-TextValue "D3" "27.31"
Set-TextValue "E3" "2.00%"

Set-TextValue "D4" "4.865"
Set-TextValue "E4" "-0.23%"

Set-TextValue "D5" "0.06429"
Set-TextValue "E5" "1.68%"

Set-TextValue "D6" "7.014"
Set-TextValue "E6" "1.33%"

Set-TextValue "D7" "1.191"
Set-TextValue "E7" "-6.77%"

Set-TextValue "D8" "0.8885"
Set-TextValue "E8" "1.78%"

Set-TextValue "D9" "0.1541"
Set-TextValue "E9" "-0.06%"

Set-TextValue "D10" "0.05177"
Set-TextValue "E10" "2.35%"

Set-TextValue "D11" "0.07503"
Set-TextValue "E11" "0.84%"

Set-TextValue "D12" "0.02900"
Set-TextValue "E12" "-2.44%"

Set-TextValue "D13" "0.08978"
Set-TextValue "E13" "-0.93%"

Set-TextValue "D14" "0.001562"
Set-TextValue "E14" "-0.55%"

Set-TextValue "D15" "0.0006368"
Set-TextValue "E15" "0.36%"

Set-TextValue "D16" "0.006148"
Set-TextValue "E16" "2.69%"

Set-TextValue "D17" "3.472"
Set-TextValue "E17" "0.52%"

Set-TextValue "D18" "3.306"
Set-TextValue "E18" "-0.50%"

Set-TextValue "D19" "2.271"
Set-TextValue "E19" "-0.05%"

Set-TextValue "E21" "0.43%"

Set-TextValue "D22" "3.918"
Set-TextValue "E22" "-0.20%"

Set-TextValue "D23" "0.1519"
Set-TextValue "E23" "10.03%"

Set-TextValue "D24" "0.04410"
Set-TextValue "E24" "0.59%"

Set-TextValue "D25" "0.001177"
Set-TextValue "E25" "1.28%"

Set-TextValue "D26" "0.003885"
Set-TextValue "E26" "-7.71%"

Set-TextValue "E28" "-1.78%"

Set-TextValue "E29" "1.62%"

Set-TextValue "D40" "0.04117"

Set-TextValue "D41" "0.006817"
Set-TextValue "E41" "-3.19%"

Set-TextValue "E42" "0.26%"

Set-TextValue "D43" "0.001911"
Set-TextValue "E43" "-15.21%"

Set-TextValue "D44" "0.01167"
Set-TextValue "E44" "3.73%"

Set-TextValue "D45" "0.00005326"
Set-TextValue "E45" "2.29%"

Set-TextValue "E46" "13.18%"
